$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 15 new blank rows before the current last two "border" rows
# (old row 27 -> new row 42, old row 28 -> new row 43), pushing everything
# after row 26 down by 15 rows.
$ws.Rows.Item(27).Resize(15).Insert()

# The inserted rows come back with no formatting, so restore the same
# per-column styling used by the existing blank rows (21-26) by copying
# formats from row 21 down onto the freshly inserted rows.
$ws.Range("A21:D21").Copy()
$ws.Range("A27:D41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the previously-blank task rows 21-24 with the new tasks.
$ws.Range("A21").Value = "20. Add review page"
$ws.Range("B21").Value = "ASAP"
$ws.Range("C21").Value = "In Progress"
$ws.Range("D21").Value = "Arthur"

$ws.Range("A22").Value = "21. Change settings address (licence, word list) to my documents"
$ws.Range("B22").Value = "Hight"
$ws.Range("C22").Value = "Open"

$ws.Range("A23").Value = "22. Add feature to calculate new news on the sites"
$ws.Range("B23").Value = "Hight"
$ws.Range("C23").Value = "Open"

$ws.Range("A24").Value = "23. Add feature to calculate new added by user"
$ws.Range("B24").Value = "Hight"
$ws.Range("C24").Value = "Open"

# Update the selection to match the author's saved cursor position.
$ws.Range("F9").Select()
